$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values for rows 10-24 ---
$ws.Range("B10").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C10").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "01/01/2018"
$ws.Range("C13").Value = "01/01/2018"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("C15").Value = "8151869 - Livia Chaguri e Carvalho"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. Comminution and solids classification2. Filtration3. Fluid mixing4. Heat exchangers5. Evaporation6. Distillation7. Absorption8. Liquid-liquid extraction"
$ws.Range("C16").Value = "1. Comminution and solids classification2. Filtration3. Fluid mixing4. Heat exchangers5. Evaporation6. Distillation7. Absorption8. Liquid-liquid extraction"
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("C18").Value = "5817372 - Simone de Fátima Medeiros Sampaio"
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aplicação de prova(s) e relatório(s)."
$ws.Range("C19").Value = "Aplicação de prova(s) e relatório(s)."
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "A média do período será definida pelo professor da disciplina. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$ws.Range("C20").Value = "A média do período será definida pelo professor da disciplina. Alunos com média final igual ou superior a 5,0 estarão aprovados, desde que tenham freqüência mínima de 70% (regimental). Alunos com média inferior a 3,0 e/ou freqüência inferior a 70% estarão reprovados (regimental). Alunos com média superior ou igual a 3,0 e inferior a 5,0 e que tenham freqüência mínima de 70% serão submetidos ao período de recuperação (regimental)."
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação."
$ws.Range("C21").Value = "A média final após a recuperação para a disciplina será a média aritmética entre a média do período e a nota da recuperação. Durante o período de recuperação, poderá ser marcada uma aula com a finalidade de sanar dúvidas e/ou revisar conceitos fundamentais. Em data posterior os alunos serão submetidos a uma prova de recuperação."
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("B23").Value = "LOQ4085 -  Operações Unitárias I  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4085 -  Operações Unitárias I  (Requisito fraco)`n"
$ws.Range("A24").ClearContents()
$ws.Range("B24").Value = "LOQ4086 -  Operações Unitárias II  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOQ4086 -  Operações Unitárias II  (Requisito fraco)`n"

# --- Update row heights for rows 13-24 ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

# --- Remove now-obsolete trailing rows 25 and 26 ---
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(25).Delete()
